$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Insert a new product row at row 7 (shifts old rows 7..10 down to 8..11) ---
$ws.Rows("7:7").Insert()

# Copy the formatting (styles + merged-cell layout) from row 6 into the newly
# inserted blank row 7 so it matches the other item rows (style ids 6/7/8/9
# and the B:G / H:K / L:M merges).
$ws.Range("A6:N6").Copy()
$ws.Range("A7:N7").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Restore the explicit row heights as they appear after the edit.
$ws.Rows("7:7").RowHeight = 25.5
$ws.Rows("9:9").RowHeight = 24.75
$ws.Rows("10:10").RowHeight = 26.25

# --- Fill in the new row 7 ("ORS 10 SACHET") ---
$ws.Range("A7").Value = 4
$ws.Range("B7").Value = "ORS 10 SACHET"
$ws.Range("H7").Value = "4:1"
$ws.Range("L7").Value = 4
$ws.Range("N7").Value = "0:0"

# --- Renumber the items that shifted down ---
$ws.Range("A8").Value = 5
$ws.Range("A9").Value = 6

# --- Update the running total (K column on the totals row, now row 10) ---
$ws.Range("K10").Value = 197
